$wb = $excel.ActiveWorkbook

# Costs and Revenues
$ws = $wb.Worksheets.Item("Costs and Revenues")
$ws.Range("B2").Value = 77388.66797673711
$ws.Range("C2").Value = 0.0
$ws.Range("D2").Value = 9992.97670278544
$ws.Range("E2").Value = 1770.0
$ws.Range("F2").Value = 35086.17445161043

# Capacities
$ws = $wb.Worksheets.Item("Capacities")
$ws.Range("C3").Value = 68.0
$ws.Range("B4").Value = 0.0
$ws.Range("D4").Value = 0.0

# PV Dispatch
$ws = $wb.Worksheets.Item("PV Dispatch")
$ws.Range("G2").Value = 13.6
$ws.Range("H2").Value = 27.2
$ws.Range("I2").Value = 34.0
$ws.Range("J2").Value = 40.8
$ws.Range("K2").Value = 47.6
$ws.Range("L2").Value = 54.4
$ws.Range("M2").Value = 61.2
$ws.Range("N2").Value = 68.0
$ws.Range("O2").Value = 61.2
$ws.Range("P2").Value = 54.4
$ws.Range("Q2").Value = 47.6
$ws.Range("R2").Value = 34.0
$ws.Range("S2").Value = 20.4
$ws.Range("T2").Value = 13.6
$ws.Range("I3").Value = 27.2
$ws.Range("J3").Value = 40.8
$ws.Range("K3").Value = 54.4
$ws.Range("L3").Value = 61.2
$ws.Range("M3").Value = 68.0
$ws.Range("N3").Value = 54.4
$ws.Range("O3").Value = 47.6
$ws.Range("P3").Value = 34.0
$ws.Range("Q3").Value = 34.0
$ws.Range("R3").Value = 20.4
$ws.Range("S3").Value = 13.6
$ws.Range("J4").Value = 6.8
$ws.Range("K4").Value = 27.2
$ws.Range("L4").Value = 47.6
$ws.Range("M4").Value = 37.58312417100181
$ws.Range("N4").Value = 54.4
$ws.Range("O4").Value = 47.6
$ws.Range("P4").Value = 27.2
$ws.Range("Q4").Value = 13.6
$ws.Range("R4").Value = 6.8

# Battery Input
$ws = $wb.Worksheets.Item("Battery Input")
$ws.Range("G2").Value = 64.3
$ws.Range("H2").Value = 14.2
$ws.Range("I2").Value = 2.8
$ws.Range("J2").Value = 1.8
$ws.Range("K2").Value = 23.53427201306108
$ws.Range("L2").Value = 33.6
$ws.Range("M2").Value = 37.8
$ws.Range("N2").Value = 42.0
$ws.Range("O2").Value = 30.0
$ws.Range("P2").Value = 150.0
$ws.Range("Q2").Value = 21.6
$ws.Range("R2").Value = 48.9
$ws.Range("T2").Value = 13.6
$ws.Range("I3").Value = 27.2
$ws.Range("J3").Value = 40.8
$ws.Range("K3").Value = 54.4
$ws.Range("L3").Value = 61.2
$ws.Range("M3").Value = 44.6
$ws.Range("N3").Value = 28.4
$ws.Range("O3").Value = 47.83079277624766
$ws.Range("P3").Value = 5.4
$ws.Range("Q3").Value = 8.0
$ws.Range("R3").Value = 20.4
$ws.Range("J4").Value = 6.8
$ws.Range("K4").Value = 27.2
$ws.Range("L4").Value = 47.6
$ws.Range("M4").Value = 14.1831241710018
$ws.Range("N4").Value = 54.4
$ws.Range("O4").Value = 47.6
$ws.Range("P4").Value = 27.2
$ws.Range("Q4").Value = 13.6
$ws.Range("R4").Value = 6.8

# State of Charge
$ws = $wb.Worksheets.Item("State of Charge")
$ws.Range("G2").Value = 183.657
$ws.Range("H2").Value = 197.715
$ws.Range("I2").Value = 200.487
$ws.Range("J2").Value = 202.269
$ws.Range("K2").Value = 225.5679292929304
$ws.Range("L2").Value = 258.8319292929305
$ws.Range("M2").Value = 296.2539292929305
$ws.Range("N2").Value = 337.8339292929305
$ws.Range("O2").Value = 367.5339292929305
$ws.Range("P2").Value = 516.0339292929305
$ws.Range("Q2").Value = 537.4179292929305
$ws.Range("R2").Value = 585.8289292929304
$ws.Range("S2").Value = 585.8289292929304
$ws.Range("I3").Value = 146.928
$ws.Range("J3").Value = 187.32
$ws.Range("K3").Value = 241.176
$ws.Range("L3").Value = 301.764
$ws.Range("M3").Value = 345.918
$ws.Range("N3").Value = 374.034
$ws.Range("O3").Value = 421.3864848484852
$ws.Range("P3").Value = 426.7324848484852
$ws.Range("Q3").Value = 434.6524848484852
$ws.Range("J4").Value = 126.732
$ws.Range("K4").Value = 153.66
$ws.Range("L4").Value = 200.784
$ws.Range("M4").Value = 214.8252929292918
$ws.Range("N4").Value = 268.6812929292918
$ws.Range("O4").Value = 315.8052929292918
$ws.Range("P4").Value = 342.7332929292918
$ws.Range("Q4").Value = 356.1972929292918

# Feed in from Type 2
$ws = $wb.Worksheets.Item("Feed in from Type 2")
$ws.Range("J2").Value = 0.0
$ws.Range("K2").Value = 1.934272013061083
$ws.Range("L2").Value = 0.0
$ws.Range("M2").Value = 0.0
$ws.Range("P2").Value = 0.0
$ws.Range("Q2").Value = 0.0
$ws.Range("I3").Value = 0.0
$ws.Range("M3").Value = 0.0
$ws.Range("N3").Value = 0.0
$ws.Range("O3").Value = 0.2307927762476609
$ws.Range("Q3").Value = 0.0
$ws.Range("M4").Value = 0.0

# Feed in from Type 3
$ws = $wb.Worksheets.Item("Feed in from Type 3")
$ws.Range("H2").Value = 0.0
$ws.Range("I2").Value = 0.0
$ws.Range("J2").Value = 0.0
$ws.Range("N2").Value = 0.0
$ws.Range("R2").Value = 48.7
$ws.Range("S2").Value = 0.0
$ws.Range("O3").Value = 0.0
$ws.Range("Q3").Value = 0.0
$ws.Range("R3").Value = 0.0
$ws.Range("J4").Value = 0.0
$ws.Range("Q4").Value = 0.0
$ws.Range("R4").Value = 0.0

# Feed in from Type 4
$ws = $wb.Worksheets.Item("Feed in from Type 4")
$ws.Range("J2").Value = 0.0
$ws.Range("N2").Value = 0.0
$ws.Range("O2").Value = 0.0
$ws.Range("P2").Value = 124.2
$ws.Range("S2").Value = 21.2
$ws.Range("T2").Value = 0.0
$ws.Range("P3").Value = 0.0
$ws.Range("S3").Value = 9.6
$ws.Range("J4").Value = 0.0
$ws.Range("L4").Value = 0.0
$ws.Range("Q4").Value = 0.0
$ws.Range("R4").Value = 0.0
